$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = [double]"1"
$ws.Range("F2").Value = [double]"0.3333333333333333"
$ws.Range("G2").Value = [double]"0.2189473333333334"
$ws.Range("H2").Value = [double]"0.656842"
$ws.Range("I2").Value = [double]"0.009402596261870986"
$ws.Range("J2").Value = [double]"0.009402596261870984"
$ws.Range("K2").Value = [double]"3"
$ws.Range("L2").Value = [double]"1"
$ws.Range("M2").Value = [double]"0.4702473333333333"
$ws.Range("N2").Value = [double]"1.410742"
$ws.Range("O2").Value = [double]"0.009034922268422819"
$ws.Range("P2").Value = [double]"0.009034922268422819"
$ws.Range("Q2").Value = [double]"0.1029593996404444"
$ws.Range("R2").Value = [double]"0.9266345967640001"
$ws.Range("S2").Value = [double]"8.495172634736734E-05"
$ws.Range("T2").Value = [double]"8.495172634736731E-05"
$ws.Range("E3").Value = [double]"1"
$ws.Range("F3").Value = [double]"0.3333333333333333"
$ws.Range("G3").Value = [double]"0.2189473333333334"
$ws.Range("H3").Value = [double]"0.656842"
$ws.Range("I3").Value = [double]"0.009402596261870986"
$ws.Range("J3").Value = [double]"0.009402596261870984"
$ws.Range("N3").Value = [double]"0.9584440000000001"
$ws.Range("O3").Value = [double]"0.006138235792679485"
$ws.Range("P3").Value = [double]"0.006138235792679485"
$ws.Range("Q3").Value = [double]"0.06994958598311111"
$ws.Range("R3").Value = [double]"0.6295462738480001"
$ws.Range("S3").Value = [double]"5.771535291873081E-05"
$ws.Range("T3").Value = [double]"5.771535291873081E-05"
$ws.Range("E4").Value = [double]"1"
$ws.Range("F4").Value = [double]"0.3333333333333333"
$ws.Range("G4").Value = [double]"0.2189473333333334"
$ws.Range("H4").Value = [double]"0.656842"
$ws.Range("I4").Value = [double]"0.009402596261870986"
$ws.Range("J4").Value = [double]"0.009402596261870984"
$ws.Range("M4").Value = [double]"1.047307"
$ws.Range("N4").Value = [double]"3.141921"
$ws.Range("O4").Value = [double]"0.02012204358311108"
$ws.Range("P4").Value = [double]"0.02012204358311108"
$ws.Range("Q4").Value = [double]"0.2293050748313334"
$ws.Range("R4").Value = [double]"2.063745673482"
$ws.Range("S4").Value = [double]"0.0001891994517757653"
$ws.Range("T4").Value = [double]"0.0001891994517757653"
$ws.Range("E5").Value = [double]"1"
$ws.Range("F5").Value = [double]"0.3333333333333333"
$ws.Range("G5").Value = [double]"0.2189473333333334"
$ws.Range("H5").Value = [double]"0.656842"
$ws.Range("I5").Value = [double]"0.009402596261870986"
$ws.Range("J5").Value = [double]"0.009402596261870984"
$ws.Range("M5").Value = [double]"50.21070966666667"
$ws.Range("N5").Value = [double]"150.632129"
$ws.Range("O5").Value = [double]"0.9647047983557866"
$ws.Range("P5").Value = [double]"0.9647047983557866"
$ws.Range("Q5").Value = [double]"10.99350098629089"
$ws.Range("R5").Value = [double]"98.94150887661802"
$ws.Range("S5").Value = [double]"0.009070729730829123"
$ws.Range("T5").Value = [double]"0.009070729730829121"
$ws.Range("I6").Value = [double]"0.8622887582286424"
$ws.Range("J6").Value = [double]"0.8622887582286423"
$ws.Range("K6").Value = [double]"3"
$ws.Range("L6").Value = [double]"1"
$ws.Range("M6").Value = [double]"0.4702473333333333"
$ws.Range("N6").Value = [double]"1.410742"
$ws.Range("O6").Value = [double]"0.009034922268422819"
$ws.Range("P6").Value = [double]"0.009034922268422819"
$ws.Range("Q6").Value = [double]"9.442150911439779"
$ws.Range("R6").Value = [double]"84.979358202958"
$ws.Range("S6").Value = [double]"0.007790711903530621"
$ws.Range("T6").Value = [double]"0.007790711903530621"
$ws.Range("I7").Value = [double]"0.8622887582286424"
$ws.Range("J7").Value = [double]"0.8622887582286423"
$ws.Range("N7").Value = [double]"0.9584440000000001"
$ws.Range("O7").Value = [double]"0.006138235792679485"
$ws.Range("P7").Value = [double]"0.006138235792679485"
$ws.Range("Q7").Value = [double]"6.414902858328445"
$ws.Range("R7").Value = [double]"57.73412572495601"
$ws.Range("S7").Value = [double]"0.0052929317193842"
$ws.Range("T7").Value = [double]"0.005292931719384199"
$ws.Range("I8").Value = [double]"0.8622887582286424"
$ws.Range("J8").Value = [double]"0.8622887582286423"
$ws.Range("M8").Value = [double]"1.047307"
$ws.Range("N8").Value = [double]"3.141921"
$ws.Range("O8").Value = [double]"0.02012204358311108"
$ws.Range("P8").Value = [double]"0.02012204358311108"
$ws.Range("Q8").Value = [double]"21.02899908971434"
$ws.Range("R8").Value = [double]"189.260991807429"
$ws.Range("S8").Value = [double]"0.01735101197430348"
$ws.Range("T8").Value = [double]"0.01735101197430347"
$ws.Range("I9").Value = [double]"0.8622887582286424"
$ws.Range("J9").Value = [double]"0.8622887582286423"
$ws.Range("M9").Value = [double]"50.21070966666667"
$ws.Range("N9").Value = [double]"150.632129"
$ws.Range("O9").Value = [double]"0.9647047983557866"
$ws.Range("P9").Value = [double]"0.9647047983557866"
$ws.Range("Q9").Value = [double]"1008.186680576225"
$ws.Range("R9").Value = [double]"9073.680125186023"
$ws.Range("S9").Value = [double]"0.831854102631424"
$ws.Range("T9").Value = [double]"0.8318541026314239"
$ws.Range("G10").Value = [double]"2.823530666666667"
$ws.Range("H10").Value = [double]"8.470592"
$ws.Range("I10").Value = [double]"0.1212552739852724"
$ws.Range("J10").Value = [double]"0.1212552739852723"
$ws.Range("K10").Value = [double]"3"
$ws.Range("L10").Value = [double]"1"
$ws.Range("M10").Value = [double]"0.4702473333333333"
$ws.Range("N10").Value = [double]"1.410742"
$ws.Range("O10").Value = [double]"0.009034922268422819"
$ws.Range("P10").Value = [double]"0.009034922268422819"
$ws.Range("Q10").Value = [double]"1.327757766584889"
$ws.Range("R10").Value = [double]"11.949819899264"
$ws.Range("S10").Value = [double]"0.001095531975093248"
$ws.Range("T10").Value = [double]"0.001095531975093247"
$ws.Range("G11").Value = [double]"2.823530666666667"
$ws.Range("H11").Value = [double]"8.470592"
$ws.Range("I11").Value = [double]"0.1212552739852724"
$ws.Range("J11").Value = [double]"0.1212552739852723"
$ws.Range("N11").Value = [double]"0.9584440000000001"
$ws.Range("O11").Value = [double]"0.006138235792679485"
$ws.Range("P11").Value = [double]"0.006138235792679485"
$ws.Range("Q11").Value = [double]"0.9020653420942222"
$ws.Range("R11").Value = [double]"8.118588078848001"
$ws.Range("S11").Value = [double]"0.0007442934628275565"
$ws.Range("T11").Value = [double]"0.0007442934628275564"
$ws.Range("G12").Value = [double]"2.823530666666667"
$ws.Range("H12").Value = [double]"8.470592"
$ws.Range("I12").Value = [double]"0.1212552739852724"
$ws.Range("J12").Value = [double]"0.1212552739852723"
$ws.Range("M12").Value = [double]"1.047307"
$ws.Range("N12").Value = [double]"3.141921"
$ws.Range("O12").Value = [double]"0.02012204358311108"
$ws.Range("P12").Value = [double]"0.02012204358311108"
$ws.Range("Q12").Value = [double]"2.957103431914667"
$ws.Range("R12").Value = [double]"26.613930887232"
$ws.Range("S12").Value = [double]"0.002439903907813726"
$ws.Range("T12").Value = [double]"0.002439903907813725"
$ws.Range("G13").Value = [double]"2.823530666666667"
$ws.Range("H13").Value = [double]"8.470592"
$ws.Range("I13").Value = [double]"0.1212552739852724"
$ws.Range("J13").Value = [double]"0.1212552739852723"
$ws.Range("M13").Value = [double]"50.21070966666667"
$ws.Range("N13").Value = [double]"150.632129"
$ws.Range("O13").Value = [double]"0.9647047983557866"
$ws.Range("P13").Value = [double]"0.9647047983557866"
$ws.Range("Q13").Value = [double]"141.7714785389298"
$ws.Range("R13").Value = [double]"1275.943306850368"
$ws.Range("S13").Value = [double]"0.1169755446395378"
$ws.Range("T13").Value = [double]"0.1169755446395378"
$ws.Range("G14").Value = [double]"0.1642436666666667"
$ws.Range("H14").Value = [double]"0.492731"
$ws.Range("I14").Value = [double]"0.007053371524214274"
$ws.Range("J14").Value = [double]"0.007053371524214274"
$ws.Range("K14").Value = [double]"3"
$ws.Range("L14").Value = [double]"1"
$ws.Range("M14").Value = [double]"0.4702473333333333"
$ws.Range("N14").Value = [double]"1.410742"
$ws.Range("O14").Value = [double]"0.009034922268422819"
$ws.Range("P14").Value = [double]"0.009034922268422819"
$ws.Range("Q14").Value = [double]"0.07723514626688889"
$ws.Range("R14").Value = [double]"0.695116316402"
$ws.Range("S14").Value = [double]"6.372666345158294E-05"
$ws.Range("T14").Value = [double]"6.372666345158294E-05"
$ws.Range("G15").Value = [double]"0.1642436666666667"
$ws.Range("H15").Value = [double]"0.492731"
$ws.Range("I15").Value = [double]"0.007053371524214274"
$ws.Range("J15").Value = [double]"0.007053371524214274"
$ws.Range("N15").Value = [double]"0.9584440000000001"
$ws.Range("O15").Value = [double]"0.006138235792679485"
$ws.Range("P15").Value = [double]"0.006138235792679485"
$ws.Range("Q15").Value = [double]"0.05247278561822222"
$ws.Range("R15").Value = [double]"0.4722550705640001"
$ws.Range("S15").Value = [double]"4.329525754899831E-05"
$ws.Range("T15").Value = [double]"4.329525754899831E-05"
$ws.Range("G16").Value = [double]"0.1642436666666667"
$ws.Range("H16").Value = [double]"0.492731"
$ws.Range("I16").Value = [double]"0.007053371524214274"
$ws.Range("J16").Value = [double]"0.007053371524214274"
$ws.Range("M16").Value = [double]"1.047307"
$ws.Range("N16").Value = [double]"3.141921"
$ws.Range("O16").Value = [double]"0.02012204358311108"
$ws.Range("P16").Value = [double]"0.02012204358311108"
$ws.Range("Q16").Value = [double]"0.1720135418056667"
$ws.Range("R16").Value = [double]"1.548121876251"
$ws.Range("S16").Value = [double]"0.0001419282492181143"
$ws.Range("T16").Value = [double]"0.0001419282492181143"
$ws.Range("G17").Value = [double]"0.1642436666666667"
$ws.Range("H17").Value = [double]"0.492731"
$ws.Range("I17").Value = [double]"0.007053371524214274"
$ws.Range("J17").Value = [double]"0.007053371524214274"
$ws.Range("M17").Value = [double]"50.21070966666667"
$ws.Range("N17").Value = [double]"150.632129"
$ws.Range("O17").Value = [double]"0.9647047983557866"
$ws.Range("P17").Value = [double]"0.9647047983557866"
$ws.Range("Q17").Value = [double]"8.24679106158878"
$ws.Range("R17").Value = [double]"74.22111955429901"
$ws.Range("S17").Value = [double]"0.006804421353995578"
$ws.Range("T17").Value = [double]"0.006804421353995578"
